# Updated with last transect data
# Appends the 10/3/2023 (serial 45202) transect sampling rows to Sheet1,
# reproducing the 11 new records (rows 77-87) added in the source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data rows: site, SAMcode, date (serial), time (fraction of day), Zoop_Code
$newRows = @(
    @{ Row = 77; Site = "RCS";  Sam = "ES1023B0600"; Date = 45202; Time = 0.39513888888888887; Zoop = 2 },
    @{ Row = 78; Site = "RD22"; Sam = "ES1023B0602"; Date = 45202; Time = 0.28888888888888892; Zoop = 3 },
    @{ Row = 79; Site = "I80";  Sam = "ES1023B0604"; Date = 45202; Time = 0.31111111111111112; Zoop = 2 },
    @{ Row = 80; Site = "LIS";  Sam = "ES1023B0605"; Date = 45202; Time = 0.37777777777777777; Zoop = 2 },
    @{ Row = 81; Site = "STTD"; Sam = "ES1023B0606"; Date = 45202; Time = 0.46041666666666670; Zoop = "N/A" },
    @{ Row = 82; Site = "BL5";  Sam = "ES1023B0607"; Date = 45202; Time = 0.39027777777777778; Zoop = 2 },
    @{ Row = 83; Site = "PRS";  Sam = "ES1023B0608"; Date = 45202; Time = 0.37361111111111112; Zoop = 3 },
    @{ Row = 84; Site = "LIB";  Sam = "ES1023B0609"; Date = 45202; Time = 0.35069444444444442; Zoop = 3 },
    @{ Row = 85; Site = "RYI";  Sam = "ES1023B0610"; Date = 45202; Time = 0.32222222222222224; Zoop = 3 },
    @{ Row = 86; Site = "RVB";  Sam = "ES1023B0611"; Date = 45202; Time = 0.29305555555555557; Zoop = 4 },
    @{ Row = 87; Site = "SHR";  Sam = "ES1023B0612"; Date = 45202; Time = 0.39444444444444443; Zoop = 2 }
)

foreach ($item in $newRows) {
    $r = $item.Row

    # Match the existing formatting used throughout the table:
    #  - column A (Site) uses the bordered/centered-vertical style seen from row 48 on
    #  - columns B-E (SAMcode/Date/Time/Zoop_Code) use the style from the earlier rows
    $ws.Range("B4:E4").Copy() | Out-Null
    $ws.Range("B$r" + ":E$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("A48").Copy() | Out-Null
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($r, 1).Value = $item.Site
    $ws.Cells.Item($r, 2).Value = $item.Sam
    $ws.Cells.Item($r, 3).Value = $item.Date
    $ws.Cells.Item($r, 4).Value = $item.Time
    $ws.Cells.Item($r, 5).Value = $item.Zoop
}

$ws.Activate() | Out-Null
$null = $ws.Range("A88").Select()
